$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0008407767361524817
$ws.Range("C2").Value = 0.00999999
$ws.Range("D2").Value = 0.00999999
$ws.Range("F2").Value = 0.00999999
$ws.Range("H2").Value = 0.00999999
$ws.Range("J2").Value = 0.002942209448842503
$ws.Range("K2").Value = 0.05453722682246739
$ws.Range("L2").Value = 0.02659832142169598
$ws.Range("N2").Value = 989.0641561766812
$ws.Range("P2").Value = 0.9302289132752822
$ws.Range("Q2").Value = 0.00999999
$ws.Range("S2").Value = 0.8578124242961173
$ws.Range("T2").Value = 0.00999999
$ws.Range("V2").Value = 0.4627176258650118
$ws.Range("W2").Value = 0.00999999
$ws.Range("Y2").Value = 0.4708379546660023
$ws.Range("Z2").Value = 0.00999999
$ws.Range("B3").Value = 0.0008698306078830272
$ws.Range("C3").Value = 0.00999999
$ws.Range("D3").Value = 0.00999999
$ws.Range("F3").Value = 0.00999999
$ws.Range("H3").Value = 0.00999999
$ws.Range("J3").Value = 0.008519318182797051
$ws.Range("K3").Value = 0.06214261988710894
$ws.Range("L3").Value = 0.01632120970154449
$ws.Range("N3").Value = 999.5946072104163
$ws.Range("P3").Value = 0.9453990189227494
$ws.Range("Q3").Value = 0.00999999
$ws.Range("S3").Value = 0.8461894449430177
$ws.Range("T3").Value = 0.00999999
$ws.Range("V3").Value = 0.390583050591964
$ws.Range("W3").Value = 0.00999999
$ws.Range("Y3").Value = 0.4742425369659101
$ws.Range("Z3").Value = 0.00999999
$ws.Range("B4").Value = 0.0005631281290117674
$ws.Range("C4").Value = 0.00999999
$ws.Range("D4").Value = 0.00999999
$ws.Range("F4").Value = 0.00999999
$ws.Range("H4").Value = 0.00999999
$ws.Range("J4").Value = 0.01311187350688162
$ws.Range("K4").Value = 0.03372627677742594
$ws.Range("L4").Value = 0.009474718929738396
$ws.Range("N4").Value = 992.3605072727597
$ws.Range("P4").Value = 0.9305693702780024
$ws.Range("Q4").Value = 0.00999999
$ws.Range("S4").Value = 0.8397960673070415
$ws.Range("T4").Value = 0.00999999
$ws.Range("V4").Value = 0.4570263727840671
$ws.Range("W4").Value = 0.00999999
$ws.Range("Y4").Value = 0.4976897318300076
$ws.Range("Z4").Value = 0.00999999
$ws.Range("B5").Value = 0.001256242409092629
$ws.Range("C5").Value = 0.00999999
$ws.Range("D5").Value = 0.00999999
$ws.Range("F5").Value = 0.00999999
$ws.Range("H5").Value = 0.00999999
$ws.Range("J5").Value = 0.03387880421067874
$ws.Range("K5").Value = 0.04972937091648252
$ws.Range("L5").Value = 0.04201619140646812
$ws.Range("N5").Value = 985.3774149017509
$ws.Range("P5").Value = 0.9277409654166023
$ws.Range("Q5").Value = 0.00999999
$ws.Range("S5").Value = 0.8274810872665389
$ws.Range("T5").Value = 0.00999999
$ws.Range("V5").Value = 0.4325976247301739
$ws.Range("W5").Value = 0.00999999
$ws.Range("Y5").Value = 0.4792050233490445
$ws.Range("Z5").Value = 0.00999999
$ws.Range("B6").Value = 0.0006142913417921537
$ws.Range("C6").Value = 0.00999999
$ws.Range("D6").Value = 0.00999999
$ws.Range("F6").Value = 0.00999999
$ws.Range("H6").Value = 0.00999999
$ws.Range("J6").Value = 0.02047658557504146
$ws.Range("K6").Value = 0.02086544593420282
$ws.Range("L6").Value = 0.02008716409916671
$ws.Range("N6").Value = 963.780305448121
$ws.Range("P6").Value = 0.9339514320391672
$ws.Range("Q6").Value = 0.00999999
$ws.Range("S6").Value = 0.8260008120408158
$ws.Range("T6").Value = 0.00999999
$ws.Range("V6").Value = 0.4440610267180243
$ws.Range("W6").Value = 0.00999999
$ws.Range("Y6").Value = 0.5806926542812648
$ws.Range("Z6").Value = 0.00999999
$ws.Range("B7").Value = 0.0007413336272968116
$ws.Range("C7").Value = 0.00999999
$ws.Range("D7").Value = 0.00999999
$ws.Range("F7").Value = 0.00999999
$ws.Range("H7").Value = 0.00999999
$ws.Range("J7").Value = 0.0161245161508551
$ws.Range("K7").Value = 0.05132643852667483
$ws.Range("L7").Value = 0.00668248218558809
$ws.Range("N7").Value = 1007.392638696611
$ws.Range("P7").Value = 0.9178902022031743
$ws.Range("Q7").Value = 0.00999999
$ws.Range("S7").Value = 0.819921609303279
$ws.Range("T7").Value = 0.00999999
$ws.Range("V7").Value = 0.3567122463085622
$ws.Range("W7").Value = 0.00999999
$ws.Range("Y7").Value = 0.5881152372711553
$ws.Range("Z7").Value = 0.00999999
$ws.Range("B8").Value = 0.001043390649704075
$ws.Range("C8").Value = 0.00999999
$ws.Range("D8").Value = 0.00999999
$ws.Range("F8").Value = 0.00999999
$ws.Range("H8").Value = 0.00999999
$ws.Range("J8").Value = 0.01961371895013996
$ws.Range("K8").Value = 0.03911845383678664
$ws.Range("L8").Value = 0.0456069965226502
$ws.Range("N8").Value = 1006.546682444541
$ws.Range("P8").Value = 0.9189024810141637
$ws.Range("Q8").Value = 0.00999999
$ws.Range("S8").Value = 0.8599605879895562
$ws.Range("T8").Value = 0.00999999
$ws.Range("V8").Value = 0.4309944859799551
$ws.Range("W8").Value = 0.00999999
$ws.Range("Y8").Value = 0.4927104146159
$ws.Range("Z8").Value = 0.00999999
$ws.Range("B9").Value = 0.0008123996085566535
$ws.Range("C9").Value = 0.00999999
$ws.Range("D9").Value = 0.00999999
$ws.Range("F9").Value = 0.00999999
$ws.Range("H9").Value = 0.00999999
$ws.Range("J9").Value = 0.005659871018233177
$ws.Range("K9").Value = 0.06091886459295771
$ws.Range("L9").Value = 0.01466130648451656
$ws.Range("N9").Value = 1001.952762502176
$ws.Range("P9").Value = 0.9127648967685222
$ws.Range("Q9").Value = 0.00999999
$ws.Range("S9").Value = 0.8324422715207481
$ws.Range("T9").Value = 0.00999999
$ws.Range("V9").Value = 0.4099450628785138
$ws.Range("W9").Value = 0.00999999
$ws.Range("Y9").Value = 0.4821355336844517
$ws.Range("Z9").Value = 0.00999999
$ws.Range("B10").Value = 0.00115202979110107
$ws.Range("C10").Value = 0.00999999
$ws.Range("D10").Value = 0.00999999
$ws.Range("F10").Value = 0.00999999
$ws.Range("H10").Value = 0.00999999
$ws.Range("J10").Value = 0.037043120563288
$ws.Range("K10").Value = 0.03351394978129035
$ws.Range("L10").Value = 0.04464602396862292
$ws.Range("N10").Value = 995.5152369105721
$ws.Range("P10").Value = 0.9260394361298609
$ws.Range("Q10").Value = 0.00999999
$ws.Range("S10").Value = 0.8479242908672455
$ws.Range("T10").Value = 0.00999999
$ws.Range("V10").Value = 0.3781451774989963
$ws.Range("W10").Value = 0.00999999
$ws.Range("Y10").Value = 0.5205832443001981
$ws.Range("Z10").Value = 0.00999999
$ws.Range("B11").Value = 0.0009377143774610566
$ws.Range("C11").Value = 0.00999999
$ws.Range("D11").Value = 0.00999999
$ws.Range("F11").Value = 0.00999999
$ws.Range("H11").Value = 0.00999999
$ws.Range("J11").Value = 0.009188135366662653
$ws.Range("K11").Value = 0.03364838690466343
$ws.Range("L11").Value = 0.0509350092463111
$ws.Range("N11").Value = 1011.893356404593
$ws.Range("P11").Value = 0.8964748208202329
$ws.Range("Q11").Value = 0.00999999
$ws.Range("S11").Value = 0.8558902279923287
$ws.Range("T11").Value = 0.00999999
$ws.Range("V11").Value = 0.373036763409751
$ws.Range("W11").Value = 0.00999999
$ws.Range("Y11").Value = 0.5098468443350546
$ws.Range("Z11").Value = 0.00999999
